$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 2, 18.77806223785356),
    @(2, 3, 9.118208071589756),
    @(2, 4, 13.73534553838044),
    @(2, 5, 14.30143610778043),
    @(2, 7, 37.65641753112325),
    @(2, 8, 16.50213586671763),
    @(2, 10, 8.622590987332289),
    @(2, 12, 12.10615681576139),
    @(2, 14, 18.38103329126546),
    @(2, 15, 26.31447673338139),
    @(3, 2, 18.37540876974513),
    @(3, 3, 8.955632715450889),
    @(3, 4, 13.74120235424087),
    @(3, 5, 14.33759399625431),
    @(3, 7, 37.68035254349111),
    @(3, 8, 16.54806836554743),
    @(3, 10, 8.63428026953096),
    @(3, 12, 12.09176867000322),
    @(3, 14, 18.43305416490738),
    @(3, 15, 26.37728779235091),
    @(4, 2, 18.12638732429583),
    @(4, 3, 8.853643172513721),
    @(4, 4, 13.74732898205941),
    @(4, 5, 14.36172259637784),
    @(4, 7, 37.70639222210352),
    @(4, 8, 16.5791071620191),
    @(4, 10, 8.641834108977175),
    @(4, 12, 12.08455082420719),
    @(4, 14, 18.46683938570733),
    @(4, 15, 26.42166251281349),
    @(5, 2, 18.02459886444669),
    @(5, 3, 8.811571893275854),
    @(5, 4, 13.75046254513692),
    @(5, 5, 14.37204013796593),
    @(5, 7, 37.7198502337235),
    @(5, 8, 16.59246827803063),
    @(5, 10, 8.645007330817153),
    @(5, 12, 12.08201849231016),
    @(5, 14, 18.4810718627603),
    @(5, 15, 26.44120301530852),
    @(6, 2, 18.00768241322013),
    @(6, 3, 8.804556263166424),
    @(6, 4, 13.75102135152952),
    @(6, 5, 14.37378265444698),
    @(6, 7, 37.72225665111556),
    @(6, 8, 16.59472989967133),
    @(6, 10, 8.64553998609555),
    @(6, 12, 12.08162277228917),
    @(6, 14, 18.48346325231411),
    @(6, 15, 26.44453562974782),
    @(7, 2, 18.12501563718129),
    @(7, 3, 8.85307780065102),
    @(7, 4, 13.74736866293991),
    @(7, 5, 14.36185977850937),
    @(7, 7, 37.70656220468469),
    @(7, 8, 16.57928447042048),
    @(7, 10, 8.64187651925228),
    @(7, 12, 12.08451501311151),
    @(7, 14, 18.4670294469006),
    @(7, 15, 26.42192014610585),
    @(8, 2, 18.63968153096922),
    @(8, 3, 9.062622761123517),
    @(8, 4, 13.73684018189271),
    @(8, 5, 14.31350344765877),
    @(8, 7, 37.66231305062341),
    @(8, 8, 16.51738442115015),
    @(8, 10, 8.626543462243658),
    @(8, 12, 12.10086183433474),
    @(8, 14, 18.3985879183034),
    @(8, 15, 26.3349266229407),
    @(9, 2, 19.62883398312021),
    @(9, 3, 9.45486086373235),
    @(9, 4, 13.73623796224681),
    @(9, 5, 14.23395942214696),
    @(9, 7, 37.665723810903),
    @(9, 8, 16.41852723169169),
    @(9, 10, 8.599450552884839),
    @(9, 12, 12.1456268877938),
    @(9, 14, 18.27896255355533),
    @(9, 15, 26.21055112679552),
    @(10, 2, 20.33582688317011),
    @(10, 3, 9.72977928981898),
    @(10, 4, 13.74795619910647),
    @(10, 5, 14.18481949378375),
    @(10, 7, 37.72333157879466),
    @(10, 8, 16.35966487987103),
    @(10, 10, 8.581341184615663),
    @(10, 12, 12.18609539587651),
    @(10, 14, 18.19990548264668),
    @(10, 15, 26.14751269976016),
    @(11, 2, 20.65172182123362),
    @(11, 3, 9.851604000841764),
    @(11, 4, 13.75591004574149),
    @(11, 5, 14.164480681913),
    @(11, 7, 37.76148847572623),
    @(11, 8, 16.33588305277899),
    @(11, 10, 8.573488883142495),
    @(11, 12, 12.20611141358297),
    @(11, 14, 18.1658452835976),
    @(11, 15, 26.12501769103586),
    @(12, 2, 20.77040970109324),
    @(12, 3, 9.897242189827983),
    @(12, 4, 13.75929719790372),
    @(12, 5, 14.15706841457145),
    @(12, 7, 37.77765120001651),
    @(12, 8, 16.32730865072119),
    @(12, 10, 8.570570602850889),
    @(12, 12, 12.21391824351289),
    @(12, 14, 18.15322033324178),
    @(12, 15, 26.1173901049682),
    @(13, 2, 20.74489141721633),
    @(13, 3, 9.887435606519226),
    @(13, 4, 13.75855106159126),
    @(13, 5, 14.15865190197806),
    @(13, 7, 37.77409415999675),
    @(13, 8, 16.32913611236419),
    @(13, 10, 8.571196654996093),
    @(13, 12, 12.21222686159106),
    @(13, 14, 18.15592721669702),
    @(13, 15, 26.11899319780495),
    @(14, 2, 20.66150569217334),
    @(14, 3, 9.855368744451985),
    @(14, 4, 13.75618120159215),
    @(14, 5, 14.16386506729855),
    @(14, 7, 37.76278388702309),
    @(14, 8, 16.33516898527388),
    @(14, 10, 8.573247689301459),
    @(14, 12, 12.20674915550038),
    @(14, 14, 18.16480115623398),
    @(14, 15, 26.12437230268772),
    @(15, 2, 20.61030460669577),
    @(15, 3, 9.835661686313783),
    @(15, 4, 13.75477839237199),
    @(15, 5, 14.16709599037062),
    @(15, 7, 37.75607898724146),
    @(15, 8, 16.33892047352653),
    @(15, 10, 8.574511190291743),
    @(15, 12, 12.2034233728624),
    @(15, 14, 18.17027221836301),
    @(15, 15, 26.12778321646644),
    @(16, 2, 20.3150575456367),
    @(16, 3, 9.721750127022375),
    @(16, 4, 13.74748901602169),
    @(16, 5, 14.18618921836721),
    @(16, 7, 37.72107803981824),
    @(16, 8, 16.36127939577237),
    @(16, 10, 8.581862091383327),
    @(16, 12, 12.18481930896453),
    @(16, 14, 18.20216962812954),
    @(16, 15, 26.14910735893606),
    @(17, 2, 20.13238583235258),
    @(17, 3, 9.651019008425092),
    @(17, 4, 13.74368773335189),
    @(17, 5, 14.19841830694765),
    @(17, 7, 37.70266411643505),
    @(17, 8, 16.37576338403268),
    @(17, 10, 8.586470247771274),
    @(17, 12, 12.1738151611223),
    @(17, 14, 18.2222245650945),
    @(17, 15, 26.1637737178136),
    @(18, 2, 20.02678575734125),
    @(18, 3, 9.610033257290654),
    @(18, 4, 13.74174837747356),
    @(18, 5, 14.20564184584129),
    @(18, 7, 37.69319821613213),
    @(18, 8, 16.38437605775907),
    @(18, 10, 8.589157060213667),
    @(18, 12, 12.16763733846726),
    @(18, 14, 18.23393882323523),
    @(18, 15, 26.17279115919172),
    @(19, 2, 19.99094353553745),
    @(19, 3, 9.596105034762429),
    @(19, 4, 13.74113423117157),
    @(19, 5, 14.20812019691654),
    @(19, 7, 37.69018661364817),
    @(19, 8, 16.38734055325778),
    @(19, 10, 8.590073014962906),
    @(19, 12, 12.16557176309961),
    @(19, 14, 18.23793586972069),
    @(19, 15, 26.1759441607301),
    @(20, 2, 20.15188745931089),
    @(20, 3, 9.658580038122686),
    @(20, 4, 13.74406683220125),
    @(20, 5, 14.19709686861617),
    @(20, 7, 37.70450787486264),
    @(20, 8, 16.37419236216628),
    @(20, 10, 8.58597594422762),
    @(20, 12, 12.17497092183986),
    @(20, 14, 18.22007114314773),
    @(20, 15, 26.16215223612207),
    @(21, 2, 20.68602432380171),
    @(21, 3, 9.864801191187226),
    @(21, 4, 13.75686712131875),
    @(21, 5, 14.162325976139),
    @(21, 7, 37.7660595336741),
    @(21, 8, 16.333385275542),
    @(21, 10, 8.57264375427285),
    @(21, 12, 12.20835195606847),
    @(21, 14, 18.16218726453684),
    @(21, 15, 26.12276814019697),
    @(22, 2, 21.0296190191764),
    @(22, 3, 9.996686076432754),
    @(22, 4, 13.76741870220813),
    @(22, 5, 14.14128910196543),
    @(22, 7, 37.81627109468639),
    @(22, 8, 16.30922943554902),
    @(22, 10, 8.564252121724142),
    @(22, 12, 12.23149075612938),
    @(22, 14, 18.12594717407129),
    @(22, 15, 26.10222116115505),
    @(23, 2, 20.84677389906266),
    @(23, 3, 9.926570326975369),
    @(23, 4, 13.76158785060979),
    @(23, 5, 14.15236249638841),
    @(23, 7, 37.78856089518354),
    @(23, 8, 16.32189165345571),
    @(23, 10, 8.56870153926236),
    @(23, 12, 12.21902146713558),
    @(23, 14, 18.14514392478753),
    @(23, 15, 26.11271179541851),
    @(24, 2, 20.14307257067903),
    @(24, 3, 9.655162695791212),
    @(24, 4, 13.7438946750451),
    @(24, 5, 14.19769369009385),
    @(24, 7, 37.70367082098621),
    @(24, 8, 16.37490173143903),
    @(24, 10, 8.586199301962298),
    @(24, 12, 12.17444793923342),
    @(24, 14, 18.22104413071792),
    @(24, 15, 26.1628834841453),
    @(25, 2, 19.36419233370138),
    @(25, 3, 9.35094926945073),
    @(25, 4, 13.73426006314632),
    @(25, 5, 14.25384366370564),
    @(25, 7, 37.65512566169307),
    @(25, 8, 16.44285521210665),
    @(25, 10, 8.606463253860813),
    @(25, 12, 12.13217268453631),
    @(25, 14, 18.30976907331365),
    @(25, 15, 26.23923180891193),
)

foreach ($item in $data) {
    $r = $item[0]
    $c = $item[1]
    $v = $item[2]
    $ws.Cells.Item($r, $c).Value = $v
}
